$d = $word.ActiveDocument

# 1. Remove the leading space before "Możliwość propozycji" so the
#    paragraph text becomes "Możliwość propozycji alternatywnej opcji lotu."
$null = $d.Content.Find.Execute(" Możliwość propozycji alternatywnej opcji", $true, $false, $false, $false, $false, $true, 1, $false, "Możliwość propozycji alternatywnej opcji", 2)

# 2. Replace "Loty mogą być dodawane tylko przez uprawnione osoby" with the
#    new first item, then insert the remaining new list items after it.
$null = $d.Content.Find.Execute(" Loty mogą być dodawane tylko przez uprawnione osoby", $true, $false, $false, $false, $false, $true, 1, $false, " Tylko administrator może tworzyć konta dla pracowników.", 2)

$p = $d.Paragraphs(35)
$p.Range.InsertParagraphAfter()
$d.Paragraphs(36).Range.Text = " Tylko pracownik może zarządzać pilotami."
$d.Paragraphs(36).Range.InsertParagraphAfter()
$d.Paragraphs(37).Range.Text = " Tylko pracownik może zarządzać samolotami."
$d.Paragraphs(37).Range.InsertParagraphAfter()
$d.Paragraphs(38).Range.Text = " Tylko pracownik może zarządzać lotami."
$d.Paragraphs(38).Range.InsertParagraphAfter()
$d.Paragraphs(39).Range.Text = " Zakaz zakupu biletu na wyprzedany lot."
$d.Paragraphs(39).Range.InsertParagraphAfter()
$d.Paragraphs(40).Range.Text = " W przypadku odwołania lotu klient powinien dostać zwrot kosztów."
$d.Paragraphs(40).Range.InsertParagraphAfter()
$d.Paragraphs(41).Range.Text = " W przypadku odwołania lotu klient powinien dostać propozycję innego lotu w tej samej lub mniejszej cenie."
$d.Paragraphs(41).Range.InsertParagraphAfter()
$d.Paragraphs(42).Range.Text = " W przypadku rezygnacji z zakupionego biletu klient dostaje zwrot 50% kosztów."

# 3. Add a period after "Zakup biletów jest wykonywany przez klientów" and
#    insert the new paragraph about automatic ticket purchase ordering.
$null = $d.Content.Find.Execute(" Zakup biletów jest wykonywany przez klientów", $true, $false, $false, $false, $false, $true, 1, $false, " Zakup biletów jest wykonywany przez klientów.", 2)

$p2 = $d.Paragraphs(43)
$p2.Range.InsertParagraphAfter()
$d.Paragraphs(44).Range.Text = " Podczas automatycznego zakupu biletów na wyprzedane loty – bilety są przydzielane klientom w odpowiedniej kolejności."
